$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
}

# Row 2
Set-TextValue "D2" "23.210.55"
Set-TextValue "E2" "  -2.89%  "

# Row 3
Set-TextValue "D3" "1.599.02"
Set-TextValue "E3" "  -3.57%  "

# Row 4
Set-TextValue "D4" "1.004"
Set-TextValue "E4" "  +0.25%  "

# Row 5
Set-TextValue "E5" "  +0.32%  "

# Row 6
Set-TextValue "D6" "301.91"
Set-TextValue "E6" "  -2.31%  "

# Row 7
Set-TextValue "E7" "  -2.98%  "

# Row 8
Set-TextValue "D8" "0.3668"
Set-TextValue "E8" "  -4.31%  "

# Row 9
Set-TextValue "D9" "48.89"
Set-TextValue "E9" "  -4.38%  "

# Row 10
Set-TextValue "E10" "  +0.27%  "

# Row 11
Set-TextValue "D11" "1.280"
Set-TextValue "E11" "  -5.50%  "

# Row 12
Set-TextValue "D12" "0.08112"
Set-TextValue "E12" "  -4.13%  "

# Row 13
Set-TextValue "D13" "22.90"
Set-TextValue "E13" "  -4.31%  "

# Row 14
Set-TextValue "D14" "6.648"
Set-TextValue "E14" "  -6.99%  "

# Row 15
Set-TextValue "D15" "7.575"
Set-TextValue "E15" "  -3.75%  "

# Row 16
Set-TextValue "E16" "  -3.16%  "

# Row 17
Set-TextValue "D17" "1.595.05"
Set-TextValue "E17" "  -3.48%  "

# Row 18
Set-TextValue "D18" "91.90"
Set-TextValue "E18" "  -3.03%  "

# Row 19
Set-TextValue "D19" "0.06825"
Set-TextValue "E19" "  -2.53%  "

# Row 20
Set-TextValue "E20" "  -6.49%  "

# Row 21
Set-TextValue "D21" "6.622"
Set-TextValue "E21" "  -4.12%  "

# Row 22
Set-TextValue "D22" "1.004"
Set-TextValue "E22" "  +0.31%  "

# Row 23
Set-TextValue "D23" "13.17"
Set-TextValue "E23" "  -3.12%  "

# Row 24
Set-TextValue "D24" "23.209.79"
Set-TextValue "E24" "  -2.93%  "

# Row 25
Set-TextValue "B25" "LidoDAOToken"
Set-TextValue "C25" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue "D25" "2.983"
Set-TextValue "E25" "  -1.95%  "

# Row 26
Set-TextValue "B26" "Toncoin"
Set-TextValue "C26" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D26" "2.362"
Set-TextValue "E26" "  -5.35%  "

# Row 27
Set-TextValue "E27" "  -4.14%  "

# Row 28
Set-TextValue "D28" "150.95"
Set-TextValue "E28" "  -1.06%  "

# Row 29
Set-TextValue "D29" "5.284"
Set-TextValue "E29" "  -3.17%  "

# Row 30
Set-TextValue "D30" "132.61"
Set-TextValue "E30" "  -4.86%  "

# Row 31
Set-TextValue "D31" "2.469"
Set-TextValue "E31" "  -1.16%  "

# Row 32
Set-TextValue "D32" "7.143"
Set-TextValue "E32" "  -8.33%  "

# Row 33
Set-TextValue "D33" "1.772.75"
Set-TextValue "E33" "  -3.29%  "

# Row 34
Set-TextValue "D34" "0.9746"
Set-TextValue "E34" "  -5.19%  "

# Row 35
Set-TextValue "D35" "0.07759"
Set-TextValue "E35" "  -3.45%  "

# Row 36
Set-TextValue "D36" "0.02788"
Set-TextValue "E36" "  -5.76%  "

# Row 37
Set-TextValue "D37" "6.288"
Set-TextValue "E37" "  -5.79%  "

# Row 38
Set-TextValue "D38" "0.2553"
Set-TextValue "E38" "  -4.82%  "

# Row 39
Set-TextValue "D39" "10.21"
Set-TextValue "E39" "  -6.90%  "

# Row 40
Set-TextValue "D40" "0.08883"
Set-TextValue "E40" "  -2.53%  "

# Row 41
Set-TextValue "D41" "1.381"
Set-TextValue "E41" "  -2.63%  "

# Row 42
Set-TextValue "E42" "  -4.73%  "

# Row 43
Set-TextValue "D43" "12.86"
Set-TextValue "E43" "  -4.64%  "

# Row 44
Set-TextValue "D44" "16.11"
Set-TextValue "E44" "  -0.91%  "

# Row 45
Set-TextValue "D45" "0.6643"
Set-TextValue "E45" "  -4.34%  "

# Row 46
Set-TextValue "D46" "2.320"
Set-TextValue "E46" "  -5.76%  "

# Row 47
Set-TextValue "E47" "  +0.33%  "

# Row 48
Set-TextValue "D48" "3.967"
Set-TextValue "E48" "  -2.60%  "

# Row 49
Set-TextValue "D49" "0.08010"
Set-TextValue "E49" "  -3.25%  "

# Row 50
Set-TextValue "D50" "131.48"
Set-TextValue "E50" "  -2.07%  "

# Row 51
Set-TextValue "D51" "1.175"
Set-TextValue "E51" "  -4.44%  "
